# -----------------------------------------------------------------------
# checkboard.xlsx — extend the letter grid with a second ("polish keys")
# block in columns N:X, mirroring the existing Q/W/E... layout, fix the
# A2 row-header value, and refresh the view/print state.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 used to be a blank "note" cell; it now carries the row value 0
# (matching the 0..9 sequence used by A3:A9).
$ws.Range("A2").Value = 0

# The new block reuses the workbook's two existing cell styles:
#   - the header/index cells use the "Uwaga" style carried by B1 (s=2)
#   - the letter cells use the "Dane wejściowe" style carried by B2 (s=1)
# Copy-pasting formats (instead of `.Style = "..."`) keeps the same
# style indexes instead of minting new ones in styles.xml.
$numStyleSource = $ws.Range("B1")
$strStyleSource = $ws.Range("B2")
$blankStyleSource = $ws.Range("E2")

$cells = @(
    @{ Cell = "O1"; Kind = "num2"; Value = 0 },
    @{ Cell = "P1"; Kind = "num2"; Value = 1 },
    @{ Cell = "Q1"; Kind = "num2"; Value = 2 },
    @{ Cell = "R1"; Kind = "num2"; Value = 3 },
    @{ Cell = "S1"; Kind = "num2"; Value = 4 },
    @{ Cell = "T1"; Kind = "num2"; Value = 5 },
    @{ Cell = "U1"; Kind = "num2"; Value = 6 },
    @{ Cell = "V1"; Kind = "num2"; Value = 7 },
    @{ Cell = "W1"; Kind = "num2"; Value = 8 },
    @{ Cell = "X1"; Kind = "num2"; Value = 9 },
    @{ Cell = "O2"; Kind = "str1"; Value = "q" },
    @{ Cell = "P2"; Kind = "str1"; Value = "w" },
    @{ Cell = "Q2"; Kind = "blank1" },
    @{ Cell = "R2"; Kind = "str1"; Value = "e" },
    @{ Cell = "S2"; Kind = "blank1" },
    @{ Cell = "T2"; Kind = "blank1" },
    @{ Cell = "U2"; Kind = "str1"; Value = "t" },
    @{ Cell = "V2"; Kind = "str1"; Value = "y" },
    @{ Cell = "W2"; Kind = "str1"; Value = "u" },
    @{ Cell = "X2"; Kind = "str1"; Value = "i" },
    @{ Cell = "N3"; Kind = "num2"; Value = 4 },
    @{ Cell = "O3"; Kind = "str1"; Value = "o" },
    @{ Cell = "P3"; Kind = "str1"; Value = "p" },
    @{ Cell = "Q3"; Kind = "str1"; Value = "a" },
    @{ Cell = "R3"; Kind = "str1"; Value = "s" },
    @{ Cell = "S3"; Kind = "str1"; Value = "d" },
    @{ Cell = "T3"; Kind = "str1"; Value = "f" },
    @{ Cell = "U3"; Kind = "str1"; Value = "g" },
    @{ Cell = "V3"; Kind = "str1"; Value = "h" },
    @{ Cell = "W3"; Kind = "str1"; Value = "j" },
    @{ Cell = "X3"; Kind = "str1"; Value = "k" },
    @{ Cell = "N4"; Kind = "num2"; Value = 2 },
    @{ Cell = "O4"; Kind = "str1"; Value = "l" },
    @{ Cell = "P4"; Kind = "str1"; Value = "z" },
    @{ Cell = "Q4"; Kind = "str1"; Value = "x" },
    @{ Cell = "R4"; Kind = "str1"; Value = "c" },
    @{ Cell = "S4"; Kind = "str1"; Value = "v" },
    @{ Cell = "T4"; Kind = "str1"; Value = "b" },
    @{ Cell = "U4"; Kind = "str1"; Value = "n" },
    @{ Cell = "V4"; Kind = "str1"; Value = "m" },
    @{ Cell = "W4"; Kind = "str1"; Value = "ą" },
    @{ Cell = "X4"; Kind = "str1"; Value = "ć" },
    @{ Cell = "N5"; Kind = "num2"; Value = 5 },
    @{ Cell = "O5"; Kind = "str1"; Value = "ę" },
    @{ Cell = "P5"; Kind = "str1"; Value = "ł" },
    @{ Cell = "Q5"; Kind = "str1"; Value = "ń" },
    @{ Cell = "R5"; Kind = "str1"; Value = "ó" },
    @{ Cell = "S5"; Kind = "str1"; Value = "ś" },
    @{ Cell = "T5"; Kind = "str1"; Value = "ż" },
    @{ Cell = "U5"; Kind = "str1"; Value = "ź" },
    @{ Cell = "V5"; Kind = "str1"; Value = "." },
    @{ Cell = "W5"; Kind = "str1"; Value = "," },
    @{ Cell = "X5"; Kind = "str1"; Value = "r" }
)

foreach ($item in $cells) {
    $target = $ws.Range($item.Cell)

    if ($item.Kind -eq "num2") {
        $numStyleSource.Copy()
        $target.PasteSpecial(-4122)
        $target.Value = $item.Value
    } elseif ($item.Kind -eq "str1") {
        $strStyleSource.Copy()
        $target.PasteSpecial(-4122)
        $target.Value = $item.Value
    } else {
        $blankStyleSource.Copy()
        $target.PasteSpecial(-4122)
    }
}

# View: the sheet had scrolled right and zoomed out to show the new block.
$excel.ActiveWindow.Zoom = 85
$ws.Range("R10").Select()

# Print setup: orientation is now explicitly portrait.
$ws.PageSetup.Orientation = 1
